# Generate Report for Handback
# - Marks the Overview status as handed back (in sync with en-US)
# - Populates "Latest Target File" / "Latest Handback File" / "Latest Handback DateTime"
#   columns for both locale sheets (zh-cn, de-de), including hyperlinks on the
#   "Latest Target File" cells (matching the existing hyperlink look used for
#   the "Source File Name" column)
# - Widens a few columns that now hold longer text

$wb = $excel.ActiveWorkbook

$hyperlinkFontColor = 15570276   # OLE/BGR encoding of RGB 6495ED (the workbook's custom HyperLink font color)

function Set-HandoffLink {
    param($ws, [string]$cellAddr, [string]$fileName, [string]$commitSha)

    $url = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/$commitSha/e2e/$fileName"
    $cell = $ws.Range($cellAddr)
    $cell.Value = $fileName
    $ws.Hyperlinks.Add($cell, $url, [System.Type]::Missing, [System.Type]::Missing, $fileName)
    $cell.Font.Underline = $true
    $cell.Font.Color = $hyperlinkFontColor
}

$commitSha = "ba9e33a27629a712daa8a033b37f7904c06f0154"

# ---------------------------------------------------------------------------
# Overview sheet: flip the handoff status text for both locale columns/rows
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("E3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F3").Value = "Handed back: in sync with en-US"

$wsOverview.Columns.Item(5).ColumnWidth = 29.17
$wsOverview.Columns.Item(6).ColumnWidth = 29.17

# ---------------------------------------------------------------------------
# zh-cn sheet: fill in target/handback file + hyperlink for both rows
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

Set-HandoffLink $wsZhCn "I2" "a95a9fa8-071b-4fb0-8712-390905c31398.md" $commitSha
$wsZhCn.Range("J2").Value = "a95a9fa8-071b-4fb0-8712-390905c31398.0e7a525808a64c4483950b9b28b6cfa38f706fa9.zh-cn.xlf"

Set-HandoffLink $wsZhCn "I3" "c0d0d805-ce1f-4c99-978b-71f547013d2b.md" $commitSha
$wsZhCn.Range("J3").Value = "c0d0d805-ce1f-4c99-978b-71f547013d2b.f60ad76955760d077e4474b316ffc4055bb8bcaf.zh-cn.xlf"

$wsZhCn.Columns.Item(3).ColumnWidth = 29.17
$wsZhCn.Columns.Item(9).ColumnWidth = 39.17
$wsZhCn.Columns.Item(10).ColumnWidth = 39.17

# ---------------------------------------------------------------------------
# de-de sheet: fill in target/handback file + hyperlink + handback datetime
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

Set-HandoffLink $wsDeDe "I2" "a95a9fa8-071b-4fb0-8712-390905c31398.md" $commitSha
$wsDeDe.Range("J2").Value = "a95a9fa8-071b-4fb0-8712-390905c31398.0e7a525808a64c4483950b9b28b6cfa38f706fa9.de-de.xlf"
$wsDeDe.Range("K2").Value = "2016-10-20 00:22:36"

Set-HandoffLink $wsDeDe "I3" "c0d0d805-ce1f-4c99-978b-71f547013d2b.md" $commitSha
$wsDeDe.Range("J3").Value = "c0d0d805-ce1f-4c99-978b-71f547013d2b.f60ad76955760d077e4474b316ffc4055bb8bcaf.de-de.xlf"
$wsDeDe.Range("K3").Value = "2016-10-20 00:22:36"

$wsDeDe.Columns.Item(3).ColumnWidth = 29.17
$wsDeDe.Columns.Item(9).ColumnWidth = 39.17
$wsDeDe.Columns.Item(10).ColumnWidth = 39.17

Write-Output "Handback report generated"
